$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper approach: set NumberFormat to Text ("@") before assigning values that
# look numeric (prices, percentages) so Excel keeps them as literal text strings
# matching the source data (inline/shared strings), then restore the Normal style
# so no extraneous per-cell formatting is left behind.
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "305.21"
Set-TextValue $ws.Range("E2") "1.49%"
Set-TextValue $ws.Range("D3") "35.85"
Set-TextValue $ws.Range("E3") "1.53%"
Set-TextValue $ws.Range("D4") "5.067"
Set-TextValue $ws.Range("E4") "0.32%"
Set-TextValue $ws.Range("D5") "0.08054"
Set-TextValue $ws.Range("E5") "1.14%"
Set-TextValue $ws.Range("D6") "1.916"
Set-TextValue $ws.Range("E6") "1.13%"
Set-TextValue $ws.Range("D7") "4.163"
Set-TextValue $ws.Range("E7") "3.19%"
Set-TextValue $ws.Range("D8") "7.836"
Set-TextValue $ws.Range("E8") "0.75%"
Set-TextValue $ws.Range("D9") "0.9301"
Set-TextValue $ws.Range("E9") "0.32%"
Set-TextValue $ws.Range("D10") "0.1332"
Set-TextValue $ws.Range("E10") "-1.23%"
Set-TextValue $ws.Range("D11") "0.1907"
Set-TextValue $ws.Range("E11") "0.57%"
Set-TextValue $ws.Range("D12") "0.09226"
Set-TextValue $ws.Range("E12") "2.00%"
Set-TextValue $ws.Range("D13") "0.03471"
Set-TextValue $ws.Range("E13") "1.30%"
Set-TextValue $ws.Range("D14") "0.09918"
Set-TextValue $ws.Range("E14") "0.06%"
Set-TextValue $ws.Range("D15") "0.001418"
Set-TextValue $ws.Range("E15") "2.47%"
Set-TextValue $ws.Range("D16") "0.006632"
Set-TextValue $ws.Range("E16") "12.29%"
Set-TextValue $ws.Range("D17") "3.613"
Set-TextValue $ws.Range("E17") "2.41%"
Set-TextValue $ws.Range("E18") "1.73%"
Set-TextValue $ws.Range("D19") "0.3423"
Set-TextValue $ws.Range("E19") "0.50%"
Set-TextValue $ws.Range("D20") "0.1336"
Set-TextValue $ws.Range("E20") "3.30%"
Set-TextValue $ws.Range("D21") "5.172"
Set-TextValue $ws.Range("E21") "2.48%"
Set-TextValue $ws.Range("E22") "5.81%"
Set-TextValue $ws.Range("D23") "0.04411"
Set-TextValue $ws.Range("E23") "-1.88%"
Set-TextValue $ws.Range("D24") "0.001238"
Set-TextValue $ws.Range("E24") "2.00%"
Set-TextValue $ws.Range("D25") "0.004706"
Set-TextValue $ws.Range("E25") "-1.30%"
Set-TextValue $ws.Range("D26") "0.0001302"
Set-TextValue $ws.Range("E26") "5.98%"
Set-TextValue $ws.Range("D27") "0.0003134"
Set-TextValue $ws.Range("E27") "4.51%"
Set-TextValue $ws.Range("D39") "0.02002"
Set-TextValue $ws.Range("E39") "6.53%"
Set-TextValue $ws.Range("D40") "0.05179"
Set-TextValue $ws.Range("E40") "8.83%"
Set-TextValue $ws.Range("D41") "0.007633"
Set-TextValue $ws.Range("D42") "0.01016"
Set-TextValue $ws.Range("E42") "-2.31%"
Set-TextValue $ws.Range("E43") "3.29%"
Set-TextValue $ws.Range("D44") "0.002104"
Set-TextValue $ws.Range("E44") "-0.20%"
Set-TextValue $ws.Range("D45") "0.01075"
Set-TextValue $ws.Range("E45") "-2.20%"
Set-TextValue $ws.Range("D46") "0.00006306"
Set-TextValue $ws.Range("E46") "0.29%"
Set-TextValue $ws.Range("D47") "0.00000000751"
Set-TextValue $ws.Range("E47") "0.21%"
Set-TextValue $ws.Range("D48") "64.96"
Set-TextValue $ws.Range("E48") "0.49%"
Set-TextValue $ws.Range("D49") "0.001603"
Set-TextValue $ws.Range("E49") "-3.39%"
Set-TextValue $ws.Range("D50") "0.00002104"
Set-TextValue $ws.Range("E50") "0.21%"
Set-TextValue $ws.Range("D51") "0.0002004"
Set-TextValue $ws.Range("E51") "0.21%"
